# chore: update Sheets via scheduled runner
# Refresh cached market-board figures (currentAveragePrice*, LevePrice*,
# LeveProfit*) for the rows whose underlying item prices moved since the
# last run. Only the affected numeric cells are touched; cells that had no
# value before (e.g. LeveProfitNQ on rows where it was blank) are populated,
# and cells that are no longer meaningful (e.g. a stale LeveProfitHQ) are
# cleared back to blank rather than left with a stale number.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H97").Value = 1609
$ws.Range("J97").Value = 0
$ws.Range("L97").Value = 0
$ws.Range("N97").Value = ""

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H128").Value = 99390
$ws.Range("J128").Value = 99390
$ws.Range("L128").Value = 99390
$ws.Range("N128").Value = -109350

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 7107.7896
$ws.Range("I132").Value = 3202.182
$ws.Range("J132").Value = 12478
$ws.Range("K132").Value = 9606.545999999998
$ws.Range("L132").Value = 37434
$ws.Range("M132").Value = -7076.545999999998
$ws.Range("N132").Value = -42494

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 67838.2
$ws.Range("I135").Value = 1213.9
$ws.Range("J135").Value = 201086.8
$ws.Range("K135").Value = 10925.1
$ws.Range("L135").Value = 1809781.2
$ws.Range("M135").Value = -8390.1
$ws.Range("N135").Value = -1814851.2

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H31").Value = 5392.6665
$ws.Range("I31").Value = 2251.875
$ws.Range("J31").Value = 30519
$ws.Range("K31").Value = 2251.875
$ws.Range("L31").Value = 30519
$ws.Range("M31").Value = -1957.875
$ws.Range("N31").Value = -31107

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 35278.1
$ws.Range("I32").Value = 20598.922
$ws.Range("K32").Value = 20598.922
$ws.Range("M32").Value = -20311.922

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3573.6667
$ws.Range("I61").Value = 3614.6428
$ws.Range("K61").Value = 3614.6428
$ws.Range("M61").Value = -3402.6428

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 9803.519
$ws.Range("I132").Value = 11228.762
$ws.Range("J132").Value = 4815.1665
$ws.Range("K132").Value = 33686.286
$ws.Range("L132").Value = 14445.4995
$ws.Range("M132").Value = -31156.286
$ws.Range("N132").Value = -19505.4995

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 3573.6667
$ws.Range("I136").Value = 3614.6428
$ws.Range("K136").Value = 10843.9284
$ws.Range("M136").Value = -8293.928400000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 4214
$ws.Range("I134").Value = 3666
$ws.Range("K134").Value = 10998
$ws.Range("M134").Value = -8463

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 298.7143
$ws.Range("I22").Value = 297.5
$ws.Range("J22").Value = 299.2
$ws.Range("K22").Value = 297.5
$ws.Range("L22").Value = 299.2
$ws.Range("M22").Value = 52.5
$ws.Range("N22").Value = -999.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 12937.444
$ws.Range("I58").Value = 50895
$ws.Range("J58").Value = 2092.4285
$ws.Range("K58").Value = 50895
$ws.Range("L58").Value = 2092.4285
$ws.Range("M58").Value = -50692
$ws.Range("N58").Value = -2498.4285

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 13818.5
$ws.Range("I99").Value = 1765
$ws.Range("J99").Value = 49979
$ws.Range("K99").Value = 1765
$ws.Range("L99").Value = 49979
$ws.Range("M99").Value = -267
$ws.Range("N99").Value = -52975

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H124").Value = 55249
$ws.Range("J124").Value = 55249
$ws.Range("L124").Value = 55249
$ws.Range("N124").Value = -60159

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H125").Value = 33343550
$ws.Range("J125").Value = 33343550
$ws.Range("L125").Value = 33343550
$ws.Range("N125").Value = -33348470

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 13818.5
$ws.Range("I126").Value = 1765
$ws.Range("J126").Value = 49979
$ws.Range("K126").Value = 5295
$ws.Range("L126").Value = 149937
$ws.Range("M126").Value = -2825
$ws.Range("N126").Value = -154877

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2670.8928
$ws.Range("I132").Value = 2510.28
$ws.Range("K132").Value = 7530.84
$ws.Range("M132").Value = -5000.84

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 2660.862
$ws.Range("I134").Value = 2303.5833
$ws.Range("K134").Value = 6910.749899999999
$ws.Range("M134").Value = -4375.749899999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 12937.444
$ws.Range("I136").Value = 50895
$ws.Range("J136").Value = 2092.4285
$ws.Range("K136").Value = 152685
$ws.Range("L136").Value = 6277.2855
$ws.Range("M136").Value = -150135
$ws.Range("N136").Value = -11377.2855

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H62").Value = 9749.25
$ws.Range("I62").Value = 7999
$ws.Range("K62").Value = 23997
$ws.Range("M62").Value = -23311

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H63").Value = 8497
$ws.Range("I63").Value = 1995
$ws.Range("K63").Value = 5985
$ws.Range("M63").Value = -5236

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H65").Value = 9749.25
$ws.Range("I65").Value = 7999
$ws.Range("K65").Value = 71991
$ws.Range("M65").Value = -68559

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H66").Value = 8497
$ws.Range("I66").Value = 1995
$ws.Range("K66").Value = 17955
$ws.Range("M66").Value = -14211

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H102").Value = 14166.333
$ws.Range("J102").Value = 14999.6
$ws.Range("L102").Value = 44998.8
$ws.Range("N102").Value = -49866.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H125").Value = 9989.5
$ws.Range("I125").Value = 9979
$ws.Range("K125").Value = 29937
$ws.Range("M125").Value = -25017

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H126").Value = 8999.666999999999
$ws.Range("I126").Value = 8999
$ws.Range("K126").Value = 26997
$ws.Range("M126").Value = -22057

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 33337950
$ws.Range("I70").Value = 33337950
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 33337950
$ws.Range("L70").Value = 0
$ws.Range("M70").Value = ""
$ws.Range("N70").Value = -33337680

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 33337950
$ws.Range("I73").Value = 33337950
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 33337950
$ws.Range("L73").Value = 0
$ws.Range("M73").Value = ""
$ws.Range("N73").Value = -33337014

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 26160.04
$ws.Range("I97").Value = 36833.89
$ws.Range("J97").Value = 2143.875
$ws.Range("K97").Value = 36833.89
$ws.Range("L97").Value = 2143.875
$ws.Range("M97").Value = -36337.89
$ws.Range("N97").Value = -3135.875

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3292.7144
$ws.Range("I132").Value = 3292.7144
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 9878.143199999999
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = ""
$ws.Range("N132").Value = -7348.143199999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2582.8333
$ws.Range("I22").Value = 2332.6667
$ws.Range("J22").Value = 2833
$ws.Range("K22").Value = 2332.6667
$ws.Range("L22").Value = 2833
$ws.Range("M22").Value = -2037.6667
$ws.Range("N22").Value = -3423

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 2582.8333
$ws.Range("I27").Value = 2332.6667
$ws.Range("J27").Value = 2833
$ws.Range("K27").Value = 2332.6667
$ws.Range("L27").Value = 2833
$ws.Range("M27").Value = -2225.6667
$ws.Range("N27").Value = -3047

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4841.6816
$ws.Range("I40").Value = 4783.5625
$ws.Range("K40").Value = 4783.5625
$ws.Range("M40").Value = -4647.5625

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2479.6667
$ws.Range("J68").Value = 748
$ws.Range("L68").Value = 748
$ws.Range("N68").Value = -2246

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 2479.6667
$ws.Range("J71").Value = 748
$ws.Range("L71").Value = 3740
$ws.Range("N71").Value = -11228

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("M2").Value = ""

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H98").Value = 30590
$ws.Range("J98").Value = 30590
$ws.Range("L98").Value = 30590
$ws.Range("N98").Value = -36580

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 6779.913
$ws.Range("I132").Value = 6953.1465
$ws.Range("K132").Value = 20859.4395
$ws.Range("M132").Value = -18329.4395
